$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes ---
# C2 was "Accuracy" -> now "Accuracy " (trailing space)
$ws.Range("C2").Value = "Accuracy "
# H3 was "-" -> now "NA"
$ws.Range("H3").Value = "NA"

# --- Alignment updates: add horizontal=left across the font-1 styled cells ---
$ws.Range("A2:F2").HorizontalAlignment = -4131
$ws.Range("G2:H2").HorizontalAlignment = -4131
$ws.Range("B3:B9").HorizontalAlignment = -4131
$ws.Range("A3:A9").HorizontalAlignment = -4131
$ws.Range("C3:H9").HorizontalAlignment = -4131

# --- Number format change for the numeric metrics block (C3:H9) ---
# 10-decimal format -> 5-decimal format
$ws.Range("C3:H9").NumberFormat = "0.00000"

# --- Selection change ---
$ws.Range("H4").Select()

# --- Page setup ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

Write-Host "Edit complete"
